$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.165.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.553.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.83%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.571.95"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.361"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.003.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.147.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.57%  "
$ws.Range("E17").Value = "  +4.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.552.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.41%  "
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("E26").Value = "  +3.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  +4.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0796"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.82%  "
$ws.Range("E31").Value = "  +2.51%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.79%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.47%  "
$ws.Range("E34").Value = "  +5.57%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "301.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("E43").Value = "  +5.84%  "
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0939"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0230"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.97%  "
